$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1792
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1792
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1792
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -1930

# Row 133
$ws.Range("H133").Value = 83899.5
$ws.Range("J133").Value = 83899.5
$ws.Range("L133").Value = 83899.5
$ws.Range("N133").Value = -94019.5

# Row 137
$ws.Range("H137").Value = 29916.457
$ws.Range("I137").Value = 1221.2222
$ws.Range("K137").Value = 3663.6666
$ws.Range("M137").Value = -1113.6666


$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4904.6924
$ws.Range("I61").Value = 3356.8572
$ws.Range("J61").Value = 6710.5
$ws.Range("K61").Value = 3356.8572
$ws.Range("L61").Value = 6710.5
$ws.Range("M61").Value = -3144.8572
$ws.Range("N61").Value = -7134.5

# Row 132
$ws.Range("H132").Value = 2403.476
$ws.Range("I132").Value = 2098.3333
$ws.Range("J132").Value = 2632.3333
$ws.Range("K132").Value = 6294.999899999999
$ws.Range("L132").Value = 7896.999899999999
$ws.Range("M132").Value = -3764.999899999999
$ws.Range("N132").Value = -12956.9999

# Row 136
$ws.Range("H136").Value = 4904.6924
$ws.Range("I136").Value = 3356.8572
$ws.Range("J136").Value = 6710.5
$ws.Range("K136").Value = 10070.5716
$ws.Range("L136").Value = 20131.5
$ws.Range("M136").Value = -7520.571599999999
$ws.Range("N136").Value = -25231.5

# Row 138
$ws.Range("H138").Value = 89000
$ws.Range("J138").Value = 89000
$ws.Range("L138").Value = 89000
$ws.Range("N138").Value = -99280


$ws = $wb.Worksheets.Item("BSM")
# Row 54
$ws.Range("H54").Value = 10000
$ws.Range("I54").Value = 10000
$ws.Range("K54").Value = 10000
$ws.Range("M54").Value = -9516

# Row 86
$ws.Range("H86").Value = 107851.16
$ws.Range("I86").Value = 2490.3845
$ws.Range("J86").Value = 336132.84
$ws.Range("K86").Value = 2490.3845
$ws.Range("L86").Value = 336132.84
$ws.Range("M86").Value = -1367.3845
$ws.Range("N86").Value = -338378.84

# Row 89
$ws.Range("H89").Value = 107851.16
$ws.Range("I89").Value = 2490.3845
$ws.Range("J89").Value = 336132.84
$ws.Range("K89").Value = 12451.9225
$ws.Range("L89").Value = 1680664.2
$ws.Range("M89").Value = -6835.922500000001
$ws.Range("N89").Value = -1691896.2


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1849.1666
$ws.Range("I31").Value = 1147.1111
$ws.Range("K31").Value = 1147.1111
$ws.Range("M31").Value = -852.1111000000001

# Row 34
$ws.Range("H34").Value = 1849.1666
$ws.Range("I34").Value = 1147.1111
$ws.Range("K34").Value = 1147.1111
$ws.Range("M34").Value = -945.1111000000001

# Row 58
$ws.Range("H58").Value = 2176244
$ws.Range("I58").Value = 3346212.2
$ws.Range("J58").Value = 3446.1428
$ws.Range("K58").Value = 3346212.2
$ws.Range("L58").Value = 3446.1428
$ws.Range("M58").Value = -3346009.2
$ws.Range("N58").Value = -3852.1428

# Row 62
$ws.Range("H62").Value = 2636.1
$ws.Range("I62").Value = 2551.375
$ws.Range("K62").Value = 2551.375
$ws.Range("M62").Value = -1927.375

# Row 65
$ws.Range("H65").Value = 2636.1
$ws.Range("I65").Value = 2551.375
$ws.Range("K65").Value = 12756.875
$ws.Range("M65").Value = -9636.875

# Row 106
$ws.Range("H106").Value = 34824.5
$ws.Range("J106").Value = 34649
$ws.Range("L106").Value = 34649
$ws.Range("N106").Value = -37173

# Row 132
$ws.Range("H132").Value = 2817.0588
$ws.Range("I132").Value = 1164.3334
$ws.Range("J132").Value = 4676.375
$ws.Range("K132").Value = 3493.0002
$ws.Range("L132").Value = 14029.125
$ws.Range("M132").Value = -963.0001999999999
$ws.Range("N132").Value = -19089.125

# Row 134
$ws.Range("H134").Value = 1351.381
$ws.Range("I134").Value = 1367.3158
$ws.Range("K134").Value = 4101.9474
$ws.Range("M134").Value = -1566.9474

# Row 136
$ws.Range("H136").Value = 2176244
$ws.Range("I136").Value = 3346212.2
$ws.Range("J136").Value = 3446.1428
$ws.Range("K136").Value = 10038636.6
$ws.Range("L136").Value = 10338.4284
$ws.Range("M136").Value = -10036086.6
$ws.Range("N136").Value = -15438.4284


$ws = $wb.Worksheets.Item("CUL")
# Row 105
$ws.Range("H105").Value = 2792.4666
$ws.Range("J105").Value = 2928.5
$ws.Range("L105").Value = 8785.5
$ws.Range("N105").Value = -14027.5

# Row 122
$ws.Range("H122").Value = 1024.5454
$ws.Range("I122").Value = 686.6667
$ws.Range("J122").Value = 1151.25
$ws.Range("K122").Value = 6180.0003
$ws.Range("L122").Value = 10361.25
$ws.Range("M122").Value = -3730.0003
$ws.Range("N122").Value = -15261.25


$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 900.75
$ws.Range("J107").Value = 1551.5
$ws.Range("L107").Value = 1551.5
$ws.Range("N107").Value = -5391.5

# Row 126
$ws.Range("H126").Value = 2573481.8
$ws.Range("I126").Value = 2780580
$ws.Range("K126").Value = 8341740
$ws.Range("M126").Value = -8339270

# Row 139
$ws.Range("H139").Value = 45765.09
$ws.Range("J139").Value = 45765.09
$ws.Range("L139").Value = 45765.09
$ws.Range("N139").Value = -56045.09


$ws = $wb.Worksheets.Item("LTW")
# Row 104
$ws.Range("H104").Value = 9239.75
$ws.Range("J104").Value = 9239.75
$ws.Range("L104").Value = 9239.75
$ws.Range("N104").Value = -16227.75

# Row 106
$ws.Range("H106").Value = 19896.5
$ws.Range("J106").Value = 19896.5
$ws.Range("L106").Value = 19896.5
$ws.Range("N106").Value = -22420.5

# Row 132
$ws.Range("H132").Value = 1661.1282
$ws.Range("I132").Value = 1342.8334
$ws.Range("K132").Value = 4028.5002
$ws.Range("M132").Value = -1498.5002

# Row 136
$ws.Range("H136").Value = 2655.3333
$ws.Range("J136").Value = 5214.2856
$ws.Range("L136").Value = 15642.8568
$ws.Range("N136").Value = -20742.8568


$ws = $wb.Worksheets.Item("WVR")
# Row 104
$ws.Range("H104").Value = 12788.333
$ws.Range("J104").Value = 12788.333
$ws.Range("L104").Value = 12788.333
$ws.Range("N104").Value = -19776.333

# Row 105
$ws.Range("H105").Value = 33968
$ws.Range("J105").Value = 33968
$ws.Range("L105").Value = 33968
$ws.Range("N105").Value = -40956

# Row 122
$ws.Range("H122").Value = 61599.77
$ws.Range("J122").Value = 2115.25
$ws.Range("L122").Value = 6345.75
$ws.Range("N122").Value = -11245.75

# Row 123
$ws.Range("H123").Value = 45306.062
$ws.Range("I123").Value = 29950
$ws.Range("J123").Value = 47499.785
$ws.Range("K123").Value = 29950
$ws.Range("L123").Value = 47499.785
$ws.Range("M123").Value = -25050
$ws.Range("N123").Value = -57299.785

# Row 132
$ws.Range("H132").Value = 1015.1852
$ws.Range("I132").Value = 861.96155
$ws.Range("K132").Value = 2585.88465
$ws.Range("M132").Value = -55.88464999999997

# Row 136
$ws.Range("H136").Value = 22225004
$ws.Range("I136").Value = 34725184
$ws.Range("J136").Value = 2464.3333
$ws.Range("K136").Value = 104175552
$ws.Range("L136").Value = 7392.999899999999
$ws.Range("M136").Value = -104173002
$ws.Range("N136").Value = -12492.9999

